$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.209.76'
$ws.Range("E2").Value = '  +11.67%  '
$ws.Range("D3").Value = '1.816.02'
$ws.Range("E3").Value = '  +8.34%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +4.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.571'
$ws.Range("E6").Value = '  +8.23%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.41'
$ws.Range("E8").Value = '  +8.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.66'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.286'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0673'
$ws.Range("E11").Value = '  +4.80%  '
$ws.Range("E12").Value = '  +3.13%  '
$ws.Range("E13").Value = '  +8.31%  '
$ws.Range("D14").Value = '1.821.10'
$ws.Range("E14").Value = '  +8.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.644'
$ws.Range("E15").Value = '  +6.65%  '
$ws.Range("D16").Value = '34.204.05'
$ws.Range("E16").Value = '  +11.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.22'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("E18").Value = '  +6.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.17'
$ws.Range("E19").Value = '  +6.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '258.04'
$ws.Range("E20").Value = '  +6.28%  '
$ws.Range("E21").Value = '  +4.54%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  +6.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.34'
$ws.Range("E24").Value = '  +2.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  +2.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.20'
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.74'
$ws.Range("E27").Value = '  +6.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.118'
$ws.Range("E28").Value = '  +4.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.04'
$ws.Range("E29").Value = '  +5.43%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.88'
$ws.Range("E31").Value = '  +12.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0524'
$ws.Range("E32").Value = '  +6.38%  '
$ws.Range("E33").Value = '  +6.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.57'
$ws.Range("E34").Value = '  +8.10%  '
$ws.Range("D35").Value = '1.528.63'
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("E37").Value = '  +6.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.633'
$ws.Range("E38").Value = '  +5.63%  '
$ws.Range("E39").Value = '  +6.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.44'
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.79'
$ws.Range("E41").Value = '  +4.43%  '
$ws.Range("E42").Value = '  +3.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.908'
$ws.Range("E43").Value = '  +8.46%  '
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0522'
$ws.Range("E45").Value = '  +4.85%  '
$ws.Range("E46").Value = '  +5.76%  '
$ws.Range("E47").Value = '  +8.65%  '
$ws.Range("E48").Value = '  +5.98%  '
$ws.Range("E49").Value = '  +17.25%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.58'
$ws.Range("E51").Value = '  +3.74%  '
